# moved test code to run_program file from file_manipulatory
#
# 1) Notes Master "datetimeFigureOut" field: 8/22/2017 -> 8/23/2017
$p = $ppt.ActivePresentation
$notesMaster = $p.NotesMaster
$dateShape = $notesMaster.Shapes.Item(2)
$dateShape.TextFrame.TextRange.Text = "8/23/2017"

# 2) Slide layout "Plot with Legend" - reposition/resize the two picture
#    placeholders (Picture Placeholder 5 / idx=12 and Picture Placeholder 7 / idx=13).
$slideMaster = $p.SlideMaster
$layout = $slideMaster.CustomLayouts.Item(1)

# Picture Placeholder 5 (idx=12): off x stays, y 987425->961547, ext 9257331x5368930 -> 9144000x5394960 EMU
$picA = $layout.Shapes.Item(4)
$picA.Top = 75.7124
$picA.Width = 720
$picA.Height = 424.80001

# Picture Placeholder 7 (idx=13): off 695325,4313238 -> 264003,4356370 ; ext 1635125x1643062 -> 2176272x1517904 EMU
$picB = $layout.Shapes.Item(5)
$picB.Left = 20.7877
$picB.Top = 343.0213
$picB.Width = 171.36
$picB.Height = 119.52001
